$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary section ---
# Valor Mora total
$ws.Range("E11").Value = 379600
# Cant. Trabajadores (number of workers included in this statement part)
$ws.Range("C13").Value = 1

# --- Replace worker table: keep a single worker (HERNANDO ALFONSO TRIVIÑO PEREZ)
#     and drop the first worker's (KAREN ALEXIS PUELLO BENITO REBOLLO) block of rows ---
$ws.Range("B16:B23").EntireRow.Delete()

# Re-order the periods for the remaining worker ascending (2403 .. 2410) and fix
# the "Valor Mora" amounts per period for the new part of the statement.
$periodos = @("2403", "2404", "2405", "2406", "2407", "2408", "2409", "2410")
$valores  = @(52000, 52000, 52000, 52000, 52000, 52000, 52000, 15600)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}

# Column D (Nombre Trabajador) shrinks now that only the shorter name remains.
$ws.Columns("D").ColumnWidth = 33.6
